$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Tätigkeit" text for the last entry (row 33, column C)
$ws.Range("C33").Value = "Kleine Änderungen bei den Tests + Versuch Bugs zu fixen"

# Update the hours value for row 33 (column B) which drives the SUM in B35
$ws.Range("B33").Value = 4

# Move the active cell selection from B34 to C34
$ws.Range("C34").Select()
